$wb = $excel.ActiveWorkbook

# --- Sheet "Productdata": scale D/F/I columns (rows 2-11) by 0.0004 ---
$ws = $wb.Worksheets.Item("Productdata")

$ws.Range("D2").Value = 0.00448
$ws.Range("F2").Value = 0.008959999999999999
$ws.Range("I2").Value = 0.08959999999999999

$ws.Range("D3").Value = 0.00496
$ws.Range("F3").Value = 0.00992
$ws.Range("I3").Value = 0.0992

$ws.Range("D4").Value = 0.004920000000000001
$ws.Range("F4").Value = 0.009840000000000002
$ws.Range("I4").Value = 0.09840000000000002

$ws.Range("D5").Value = 0.00444
$ws.Range("F5").Value = 0.008880000000000001
$ws.Range("I5").Value = 0.0888

$ws.Range("D6").Value = 0.00048
$ws.Range("F6").Value = 0.00096
$ws.Range("I6").Value = 0.009600000000000001

$ws.Range("D7").Value = 0.00048
$ws.Range("F7").Value = 0.00096
$ws.Range("I7").Value = 0.009600000000000001

$ws.Range("D8").Value = 0.0004400000000000001
$ws.Range("F8").Value = 0.0008800000000000001
$ws.Range("I8").Value = 0.008800000000000002

$ws.Range("D9").Value = 0.00004
$ws.Range("F9").Value = 0.00008000000000000001
$ws.Range("I9").Value = 0.0008

$ws.Range("D10").Value = 0.00004
$ws.Range("F10").Value = 0.00008000000000000001
$ws.Range("I10").Value = 0.0008

$ws.Range("D11").Value = 0.00004
$ws.Range("F11").Value = 0.00008000000000000001
$ws.Range("I11").Value = 0.0008

# --- Sheet "ForcastedStandardDeviation": zero out B9:E11 (safety stock rows now averaged to 0) ---
$ws2 = $wb.Worksheets.Item("ForcastedStandardDeviation")

$ws2.Range("B9:E11").Value = 0
